# Update workbook to reflect data through 2022-12-03
# (diff shows 11-24 -> 11-25 labels plus updated Dec/Total figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) and the workbook's sheet reference
$ws.Name = "Through 2022-11-25"

# Update the header label in I1 ("2022 (through 11-24)" -> "2022 (through 11-25)")
$ws.Range("I1").Value = "2022 (through 11-25)"

# Update December value (row 12) and Total value (row 14) in column I
$ws.Range("I12").Value = 92
$ws.Range("I14").Value = 1490
